$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -0.2503095220229778
$ws.Range("E2").Value = -0.1763863393679236
$ws.Range("F2").Value = -0.7965243440797991
$ws.Range("G2").Value = 0.0636078949626635
$ws.Range("H2").Value = 0.15052815998841365
$ws.Range("I2").Value = 42.24
$ws.Range("J2").Value = 0.010599999999999454

$ws.Range("D3").Value = 0.0446605239657579
$ws.Range("E3").Value = -0.0031689563581591175
$ws.Range("F3").Value = -0.8844470535815837
$ws.Range("G3").Value = 0.020571135597453636
$ws.Range("H3").Value = 0.08582523873282964
$ws.Range("I3").Value = 37.38
$ws.Range("J3").Value = 0.014599999999999369

$ws.Range("D4").Value = -0.027543300257815494
$ws.Range("E4").Value = 0.005963237582247973
$ws.Range("F4").Value = -0.9177263510549409
$ws.Range("G4").Value = 0.00975936843675343
$ws.Range("H4").Value = 0.05523987695624826
$ws.Range("I4").Value = 35.88
$ws.Range("J4").Value = 0.031399999999999866

$ws.Range("D5").Value = -0.0013303825928988084
$ws.Range("E5").Value = 0.01252090016666907
$ws.Range("F5").Value = -0.93359115188267
$ws.Range("G5").Value = 0.006168462354453858
$ws.Range("H5").Value = 0.04235813352471028
$ws.Range("I5").Value = 31.66
$ws.Range("J5").Value = 0.033800000000000094

$ws.Range("D6").Value = -0.048219976174124865
$ws.Range("E6").Value = -0.0499739180099125
$ws.Range("F6").Value = -0.941345699488123
$ws.Range("G6").Value = 0.003739485570970672
$ws.Range("H6").Value = 0.017471802426589172
$ws.Range("I6").Value = 34.84
$ws.Range("J6").Value = 0.12240000000000023
